# The deck's Design/theme is switched from the "Integral" (Red Violet)
# colour palette to the built-in "Office Theme" colour palette.
#
# PowerPoint stores the 12 theme colour slots (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) on the active design's Theme object, in
# that exact order, reachable from the slide master. Re-pointing every
# slot to the Office Theme's RGB values reproduces the re-colouring the
# Design Gallery performs when a new theme is applied.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

function HexToRgbVal($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index -> hex, in ThemeColorScheme order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink (Office Theme built-in palette).
$officeThemeColors = @("000000", "FFFFFF", "44546A", "E7E6E6", "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47", "0563C1", "954F72")

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToRgbVal $officeThemeColors[$i - 1]
}
